$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 283; Excel shifts rows 283:403 down to 284:404
# and carries formatting down from the row below (matches the style `s="2"`
# already present on the date column).
$ws.Rows.Item(283).Insert()

# Populate the newly inserted row 283 with the new data record.
$ws.Cells.Item(283, 1).Value = 3
$ws.Cells.Item(283, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(283, 3).Value = "Coquimbo"
$ws.Cells.Item(283, 4).Value = 44784
$ws.Cells.Item(283, 5).Value = 5
$ws.Cells.Item(283, 6).Value = 100112040
$ws.Cells.Item(283, 7).Value = "Cilantro"
$ws.Cells.Item(283, 8).Value = "Sin especificar"
$ws.Cells.Item(283, 9).Value = "Primera"
$ws.Cells.Item(283, 10).Value = 225
$ws.Cells.Item(283, 11).Value = 4000
$ws.Cells.Item(283, 12).Value = 4500
$ws.Cells.Item(283, 13).Value = 4244
$ws.Cells.Item(283, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(283, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(283, 16).Value = 1415
$ws.Cells.Item(283, 17).Value = 3
$ws.Cells.Item(283, 18).Value = "Hortaliza"
